$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$th = $sm.Theme
$tv = $th.ThemeVariants
Write-Host "Count: $($tv.Count)"
for ($i = 1; $i -le $tv.Count; $i++) {
    $v = $tv.Item($i)
    Write-Host "Variant $i : Name=$($v.Name) Id=$($v.Id)"
}
